$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.422.59"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "3.327.54"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'586.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.09%  "
$ws.Range("D6").Value = "'183.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").Value = "'0.648"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.23%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -2.13%  "
$ws.Range("D10").Value = "'6.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.91%  "
$ws.Range("D11").Value = "'0.403"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").Value = "3.902.76"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("E13").Value = "  -4.40%  "
$ws.Range("D14").Value = "66.452.40"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").Value = "'26.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.07%  "
$ws.Range("D16").Value = "3.328.40"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("E17").Value = "  -2.08%  "
$ws.Range("D18").Value = "'431.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("D19").Value = "'13.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.61%  "
$ws.Range("E20").Value = "  -2.83%  "
$ws.Range("D21").Value = "'7.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.92%  "
$ws.Range("D22").Value = "'72.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.04%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E24").Value = "  +0.66%  "
$ws.Range("D25").Value = "3.452.26"
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("D26").Value = "'0.517"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("E27").Value = "  +3.57%  "
$ws.Range("D28").Value = "'0.0000115"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.39%  "
$ws.Range("D29").Value = "'9.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("D32").Value = "'22.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.90%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("D34").Value = "'5.23"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.14%  "
$ws.Range("E35").Value = "  -3.11%  "
$ws.Range("D36").Value = "'6.63"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.28%  "
$ws.Range("D37").Value = "'159.73"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("E38").Value = "  -2.78%  "
$ws.Range("D39").Value = "'1.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.26%  "
$ws.Range("D40").Value = "2.889.36"
$ws.Range("E40").Value = "  +1.91%  "
$ws.Range("D41").Value = "'26.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.95%  "
$ws.Range("D42").Value = "'0.769"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.92%  "
$ws.Range("D43").Value = "'4.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.78%  "
$ws.Range("D44").Value = "'40.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("D46").Value = "'6.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.15%  "
$ws.Range("E47").Value = "  -1.90%  "
$ws.Range("D48").Value = "'23.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.58%  "
$ws.Range("D49").Value = "'318.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.66%  "
$ws.Range("D50").Value = "'0.0273"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("E51").Value = "  +4.98%  "
